$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.057.67"
$ws.Range("E2").Value = "  +2.65%  "

$ws.Range("D3").Value = "2.697.56"
$ws.Range("E3").Value = "  +2.02%  "

$ws.Range("E4").Value = "  -0.06%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "522.79"
$ws.Range("E5").Value = "  +0.91%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "148.90"
$ws.Range("E6").Value = "  +1.16%  "

$ws.Range("E7").Value = "  +0.10%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.580"
$ws.Range("E8").Value = "  +1.01%  "

$ws.Range("D9").Value = "2.718.63"
$ws.Range("E9").Value = "  +1.72%  "

$ws.Range("E10").Value = "  -0.24%  "

$ws.Range("E11").Value = "  -0.11%  "

$ws.Range("E12").Value = "  +0.63%  "

$ws.Range("E13").Value = "  +1.25%  "

$ws.Range("D14").Value = "3.169.08"
$ws.Range("E14").Value = "  +2.14%  "

$ws.Range("D15").Value = "61.040.46"
$ws.Range("E15").Value = "  +2.74%  "

$ws.Range("B16").Value = "Avalanche"
$ws.Range("C16").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "21.53"
$ws.Range("E16").Value = "  +1.09%  "

$ws.Range("B17").Value = "WrappedEther"
$ws.Range("C17").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D17").Value = "2.826.29"
$ws.Range("E17").Value = "  +6.06%  "

$ws.Range("E18").Value = "  +0.72%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "355.63"
$ws.Range("E19").Value = "  +2.45%  "

$ws.Range("E20").Value = "  -0.65%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.61"
$ws.Range("E21").Value = "  +0.76%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.36"
$ws.Range("E22").Value = "  +2.86%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.00"
$ws.Range("E23").Value = "  -0.06%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "62.97"
$ws.Range("E24").Value = "  +2.28%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.425"

$ws.Range("E26").Value = "  +4.17%  "

$ws.Range("E27").Value = "  -0.35%  "

$ws.Range("D28").Value = "0.0₃0830"
$ws.Range("E28").Value = "  +1.16%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.30"
$ws.Range("E29").Value = "  +1.38%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.79"
$ws.Range("E30").Value = "  +4.21%  "

$ws.Range("E31").Value = "  +0.13%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "19.22"
$ws.Range("E32").Value = "  +0.48%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "150.19"
$ws.Range("E34").Value = "  +0.22%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.22"
$ws.Range("E35").Value = "  +3.59%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.955"
$ws.Range("E36").Value = "  -8.87%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.23"
$ws.Range("E37").Value = "  +4.58%  "

$ws.Range("E38").Value = "  +10.01%  "

$ws.Range("E39").Value = "  +1.33%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "36.82"
$ws.Range("E40").Value = "  +0.15%  "

$ws.Range("E41").Value = "  +0.22%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "286.73"
$ws.Range("E42").Value = "  -0.41%  "

$ws.Range("B43").Value = "Mantle"
$ws.Range("C43").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.618"
$ws.Range("E43").Value = "  -0.18%  "

$ws.Range("B44").Value = "Stellar"
$ws.Range("C44").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0994"
$ws.Range("E44").Value = "  +0.36%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "20.08"
$ws.Range("E45").Value = "  +1.61%  "

$ws.Range("D46").Value = "2.146.26"
$ws.Range("E46").Value = "  +7.66%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.994"
$ws.Range("E47").Value = "  +0.11%  "

$ws.Range("E48").Value = "  -0.26%  "

$ws.Range("E49").Value = "  +3.53%  "

$ws.Range("E50").Value = "  +0.85%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "19.30"
$ws.Range("E51").Value = "  +3.17%  "
